$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for tag "1065" right after the current row 2 (tag 1037)
$ws.Range("A3").EntireRow.Insert()

# Insert a new row for tag "1404" right before the row that (after the
# previous insert) holds tag "1406". Before this insert that is row 8.
$ws.Range("A8").EntireRow.Insert()

# --- Populate column A (tag) as text, matching the style of the existing
#     tag cells (bold, bordered, centered) ---
function Set-TagCell($addr, $value) {
    $ws.Range("A2").Copy($ws.Range($addr))
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

Set-TagCell "A2" "1037"
Set-TagCell "A3" "1065"
Set-TagCell "A4" "1172"
Set-TagCell "A5" "1370"
Set-TagCell "A6" "1392"
Set-TagCell "A7" "1399"
Set-TagCell "A8" "1404"
Set-TagCell "A9" "1406"
Set-TagCell "A10" "1412"

# --- Populate columns B:G with numeric data ---
$data = @(
    @(0, 7.14, 13.10423286807256, 2.025939589265073, 3.37347590611476, 28),
    @(0, 7.13, 22.29447150050382, 3.626556528030944, 14.58340096931597, 86),
    @(0, 7.12, 11.43705331809856, 1.920757794620315, 4.20277796320079, 19),
    @(0, 7.03, 10.64927355675065, 0.7650142517193413, 7.006603063537734, 0),
    @(0, 7.26, 16.02409046434199, 3.174723496625334, 6.040278558458219, 55),
    @(0, 7.28, 8.113389967505748, 1.396139937127157, 4.685676507111217, 6),
    @(0, 7.18, 19.11824566937721, 1.6266717308582, 13.57207349789775, 9),
    @(0, 7.07, 16.03109597208861, 2.206187332209771, 7.832334238509944, 33),
    @(0, 7.2, 18.04385261486464, 3.155559826166355, 6.328767076480719, 56)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
